# Fix property_category values that were incorrectly left as "land" in the
# 建物 (building) and 汽車 (car) sheets of the property-declaration workbook.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: column I = property_category, data rows 2-12
$wsBuilding = $wb.Worksheets.Item("建物")
for ($row = 2; $row -le 12; $row++) {
    $wsBuilding.Cells.Item($row, 9).Value = "building"
}

# 汽車 (car) sheet: column H = property_category, data rows 2-3
$wsCar = $wb.Worksheets.Item("汽車")
for ($row = 2; $row -le 3; $row++) {
    $wsCar.Cells.Item($row, 8).Value = "car"
}
